$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 7795
$ws1.Range("F5").Value = 55
$ws1.Range("F6").Value = 569
$ws1.Range("F7").Value = 1184
$ws1.Range("F8").Value = 210
$ws1.Range("F9").Value = 23
$ws1.Range("F10").Value = 172

# Sheet "演出" (show)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 14

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 7795
$ws4.Range("F5").Value = 55
$ws4.Range("F6").Value = 569
$ws4.Range("F7").Value = 1184
$ws4.Range("F8").Value = 210
$ws4.Range("F9").Value = 14
$ws4.Range("F10").Value = 23
$ws4.Range("F11").Value = 172
